$d = $word.ActiveDocument

# The first paragraph currently holds a single run:
#   "This is a Microsoft word document."
# We need to turn that into four runs:
#   "This is a Microsoft word document."
#   " ("
#   "Changed main"
#   ")"
# Range.InsertXML *replaces* the contents of the range it's called on, so we
# target exactly paragraph 1's Range (which covers just its existing text,
# not the paragraph mark) and replace it with the full four-run OOXML,
# preserving the paragraph's original identity attributes.

$target = $d.Paragraphs(1).Range

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"
                  xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
        <w:body>
          <w:p w14:paraId="5ADF5830" w14:textId="42E3A3E7" w:rsidR="00384372" w:rsidRDefault="00094D0B">
            <w:r><w:t>This is a Microsoft word document.</w:t></w:r>
            <w:r><w:t xml:space="preserve"> (</w:t></w:r>
            <w:r><w:t>Changed main</w:t></w:r>
            <w:r><w:t>)</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

[void]$target.InsertXML($xml)
